$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently holds OCR'd Python source split across columns A:F,
# rows 1-9. We are turning that raw text grid into a pandas-DataFrame-style
# export: add a numeric column header row (0..5) above, and a numeric row
# index column (0..8) to the left, shifting the original content down one
# row and right one column.

# Shift existing data down by one row, then right by one column.
$ws.Rows.Item(1).Insert()
$ws.Columns.Item(1).Insert()

# New header row: column indices 0..5 across B1:G1.
$headerVals = @(0, 1, 2, 3, 4, 5)
for ($i = 0; $i -lt $headerVals.Length; $i++) {
    $ws.Cells.Item(1, 2 + $i).Value = $headerVals[$i]
}

# New index column: row indices 0..8 down A2:A10.
for ($i = 0; $i -le 8; $i++) {
    $ws.Cells.Item(2 + $i, 1).Value = $i
}

# Style the header row + index column like a DataFrame export: bold text,
# thin box border around each cell, centered horizontally and top-aligned
# vertically.
#
# Build the format once on a scratch cell, then fan it out with
# Copy/PasteSpecial(xlPasteFormats) so every target cell picks up the same
# single finished style instead of each Range.Property= assignment minting
# its own intermediate style record.
$scratch = $ws.Range("Z100")
$scratch.Font.Bold = $true
$scratch.Borders.LineStyle = 1
$scratch.HorizontalAlignment = -4108
$scratch.VerticalAlignment = -4160

$scratch.Copy()
$headerRange = $ws.Range("B1:G1")
$indexRange = $ws.Range("A2:A10")
$headerRange.PasteSpecial(-4122)
$indexRange.PasteSpecial(-4122)

$scratch.Clear()

Write-Host "Converted OCR text grid into DataFrame-style layout"
